$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 2609.5386
$ws.Cells.Item(113, 9).Value = 2642.5
$ws.Cells.Item(113, 10).Value = 2499.6667
$ws.Cells.Item(113, 11).Value = 2642.5
$ws.Cells.Item(113, 12).Value = 2499.6667
$ws.Cells.Item(113, 13).Value = 611.5
$ws.Cells.Item(113, 14).Value = -9007.6667
$ws.Cells.Item(137, 8).Value = 1244.8235
$ws.Cells.Item(137, 9).Value = 1005.0833
$ws.Cells.Item(137, 11).Value = 3015.2499
$ws.Cells.Item(137, 13).Value = -465.2498999999998
$ws.Cells.Item(138, 8).Value = 2468.7256
$ws.Cells.Item(138, 9).Value = 1455.9546
$ws.Cells.Item(138, 10).Value = 3237.0344
$ws.Cells.Item(138, 11).Value = 4367.8638
$ws.Cells.Item(138, 12).Value = 9711.1032
$ws.Cells.Item(138, 13).Value = 772.1361999999999
$ws.Cells.Item(138, 14).Value = -19991.1032

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 345972.84
$ws.Cells.Item(32, 9).Value = 3182.0532
$ws.Cells.Item(32, 10).Value = 3202562.8
$ws.Cells.Item(32, 11).Value = 3182.0532
$ws.Cells.Item(32, 12).Value = 3202562.8
$ws.Cells.Item(32, 13).Value = -2895.0532
$ws.Cells.Item(32, 14).Value = -3203136.8
$ws.Cells.Item(97, 8).Value = 700.62164
$ws.Cells.Item(97, 9).Value = 582.5357
$ws.Cells.Item(97, 10).Value = 1068
$ws.Cells.Item(97, 11).Value = 582.5357
$ws.Cells.Item(97, 12).Value = 1068
$ws.Cells.Item(97, 13).Value = -86.53570000000002
$ws.Cells.Item(97, 14).Value = -2060
$ws.Cells.Item(102, 8).Value = 2254.9644
$ws.Cells.Item(102, 9).Value = 631.9
$ws.Cells.Item(102, 10).Value = 3156.6667
$ws.Cells.Item(102, 11).Value = 631.9
$ws.Cells.Item(102, 12).Value = 3156.6667
$ws.Cells.Item(102, 13).Value = 990.1
$ws.Cells.Item(102, 14).Value = -6400.6667
$ws.Cells.Item(122, 8).Value = 61559.63
$ws.Cells.Item(122, 9).Value = 30774.666
$ws.Cells.Item(122, 10).Value = 169307
$ws.Cells.Item(122, 11).Value = 92323.99800000001
$ws.Cells.Item(122, 12).Value = 507921
$ws.Cells.Item(122, 13).Value = -89873.99800000001
$ws.Cells.Item(122, 14).Value = -512821
$ws.Cells.Item(125, 8).Value = 56740
$ws.Cells.Item(125, 10).Value = 56740
$ws.Cells.Item(125, 12).Value = 56740
$ws.Cells.Item(125, 14).Value = -66580
$ws.Cells.Item(132, 8).Value = 22246060
$ws.Cells.Item(132, 10).Value = 85621.836
$ws.Cells.Item(132, 12).Value = 256865.508
$ws.Cells.Item(132, 14).Value = -261925.508

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(24, 8).Value = 1195
$ws.Cells.Item(24, 9).Value = 434
$ws.Cells.Item(24, 11).Value = 434
$ws.Cells.Item(24, 13).Value = -199
$ws.Cells.Item(99, 8).Value = 1773.1666
$ws.Cells.Item(99, 9).Value = 1181.3636
$ws.Cells.Item(99, 10).Value = 2703.1428
$ws.Cells.Item(99, 11).Value = 1181.3636
$ws.Cells.Item(99, 12).Value = 2703.1428
$ws.Cells.Item(99, 13).Value = 316.6364000000001
$ws.Cells.Item(99, 14).Value = -5699.1428

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1100.9565
$ws.Cells.Item(16, 9).Value = 1001
$ws.Cells.Item(16, 10).Value = 1210
$ws.Cells.Item(16, 11).Value = 1001
$ws.Cells.Item(16, 12).Value = 1210
$ws.Cells.Item(16, 13).Value = -714
$ws.Cells.Item(16, 14).Value = -1784
$ws.Cells.Item(31, 8).Value = 3532.2307
$ws.Cells.Item(31, 9).Value = 3254.9333
$ws.Cells.Item(31, 10).Value = 3910.3635
$ws.Cells.Item(31, 11).Value = 3254.9333
$ws.Cells.Item(31, 12).Value = 3910.3635
$ws.Cells.Item(31, 13).Value = -2959.9333
$ws.Cells.Item(31, 14).Value = -4500.363499999999
$ws.Cells.Item(34, 8).Value = 3532.2307
$ws.Cells.Item(34, 9).Value = 3254.9333
$ws.Cells.Item(34, 10).Value = 3910.3635
$ws.Cells.Item(34, 11).Value = 3254.9333
$ws.Cells.Item(34, 12).Value = 3910.3635
$ws.Cells.Item(34, 13).Value = -3052.9333
$ws.Cells.Item(34, 14).Value = -4314.363499999999
$ws.Cells.Item(62, 8).Value = 1600
$ws.Cells.Item(62, 9).Value = 1400
$ws.Cells.Item(62, 10).Value = 2000
$ws.Cells.Item(62, 11).Value = 1400
$ws.Cells.Item(62, 12).Value = 2000
$ws.Cells.Item(62, 13).Value = -776
$ws.Cells.Item(62, 14).Value = -3248
$ws.Cells.Item(65, 8).Value = 1600
$ws.Cells.Item(65, 9).Value = 1400
$ws.Cells.Item(65, 10).Value = 2000
$ws.Cells.Item(65, 11).Value = 7000
$ws.Cells.Item(65, 12).Value = 10000
$ws.Cells.Item(65, 13).Value = -3880
$ws.Cells.Item(65, 14).Value = -16240
$ws.Cells.Item(107, 8).Value = 550.8
$ws.Cells.Item(107, 9).Value = 444.72726
$ws.Cells.Item(107, 10).Value = 730.3077
$ws.Cells.Item(107, 11).Value = 444.72726
$ws.Cells.Item(107, 12).Value = 730.3077
$ws.Cells.Item(107, 13).Value = 1475.27274
$ws.Cells.Item(107, 14).Value = -4570.3077
$ws.Cells.Item(113, 8).Value = 1100.9565
$ws.Cells.Item(113, 9).Value = 1001
$ws.Cells.Item(113, 10).Value = 1210
$ws.Cells.Item(113, 11).Value = 1001
$ws.Cells.Item(113, 12).Value = 1210
$ws.Cells.Item(113, 13).Value = 1169
$ws.Cells.Item(113, 14).Value = -5550
$ws.Cells.Item(122, 8).Value = 1152.9166
$ws.Cells.Item(122, 9).Value = 1029.75
$ws.Cells.Item(122, 10).Value = 1276.0834
$ws.Cells.Item(122, 11).Value = 3089.25
$ws.Cells.Item(122, 12).Value = 3828.2502
$ws.Cells.Item(122, 13).Value = -639.25
$ws.Cells.Item(122, 14).Value = -8728.2502
$ws.Cells.Item(132, 8).Value = 79863.16
$ws.Cells.Item(132, 9).Value = 2030
$ws.Cells.Item(132, 11).Value = 6090
$ws.Cells.Item(132, 13).Value = -3560

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1057.0714
$ws.Cells.Item(5, 9).Value = 1057.0714
$ws.Cells.Item(5, 11).Value = 3171.2142
$ws.Cells.Item(5, 13).Value = -3059.2142
$ws.Cells.Item(122, 8).Value = 785.73334
$ws.Cells.Item(122, 9).Value = 310
$ws.Cells.Item(122, 11).Value = 2790
$ws.Cells.Item(122, 13).Value = -340
$ws.Cells.Item(135, 8).Value = 1057.0714
$ws.Cells.Item(135, 9).Value = 1057.0714
$ws.Cells.Item(135, 11).Value = 9513.642600000001
$ws.Cells.Item(135, 13).Value = -6978.642600000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2005.8064
$ws.Cells.Item(68, 9).Value = 1965.7142
$ws.Cells.Item(68, 10).Value = 2090
$ws.Cells.Item(68, 11).Value = 1965.7142
$ws.Cells.Item(68, 12).Value = 2090
$ws.Cells.Item(68, 13).Value = -1216.7142
$ws.Cells.Item(68, 14).Value = -3588
$ws.Cells.Item(71, 8).Value = 2005.8064
$ws.Cells.Item(71, 9).Value = 1965.7142
$ws.Cells.Item(71, 10).Value = 2090
$ws.Cells.Item(71, 11).Value = 9828.571
$ws.Cells.Item(71, 12).Value = 10450
$ws.Cells.Item(71, 13).Value = -6084.571
$ws.Cells.Item(71, 14).Value = -17938

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 5168.6665
$ws.Cells.Item(15, 9).Value = 3750
$ws.Cells.Item(15, 10).Value = 8006
$ws.Cells.Item(15, 11).Value = 3750
$ws.Cells.Item(15, 12).Value = 8006
$ws.Cells.Item(15, 13).Value = -3462
$ws.Cells.Item(15, 14).Value = -8582
$ws.Cells.Item(96, 8).Value = 73035
$ws.Cells.Item(96, 9).Value = 1425
$ws.Cells.Item(96, 10).Value = 84970
$ws.Cells.Item(96, 11).Value = 1425
$ws.Cells.Item(96, 12).Value = 84970
$ws.Cells.Item(96, 13).Value = -52
$ws.Cells.Item(96, 14).Value = -87716
